$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Radius of the Earth" scratch calculation block
# (rows 18, 20, 22 under columns C:G) - delete the rows entirely.
$ws.Range("A18:A22").EntireRow.Delete()

# Add a new column F with two merged/labelled blocks that flag the two
# data tables as "Data to be scrubbed" (top block, rows 6-9) and
# "Corrected Data" (bottom block, rows 11-14). The bottom block is
# written first so "Corrected Data" lands before "Data to be scrubbed"
# in the shared-strings table (matches the saved workbook's order).
$bottom = $ws.Range("F11:F14")
$bottom.Merge()
$bottom.Value = "Corrected Data"
$bottom.Font.Name = "Calibri"
$bottom.Font.Size = 12
$bottom.Font.Bold = $true
$bottom.Font.Italic = $true
$bottom.Font.Color = 16711680
$bottom.Interior.ThemeColor = 9
$bottom.Interior.TintAndShade = 0.79998168889431442
$bottom.HorizontalAlignment = -4108
$bottom.VerticalAlignment = -4108

$top = $ws.Range("F6:F9")
$top.Merge()
$top.Value = "Data to be scrubbed"
$top.Font.Name = "Calibri"
$top.Font.Size = 11
$top.Font.Bold = $true
$top.Font.Italic = $true
$top.Font.Color = 16711680
$top.Interior.ThemeColor = 7
$top.Interior.TintAndShade = 0.79998168889431442
$top.HorizontalAlignment = -4108
$top.VerticalAlignment = -4108

# New column F is a bit wider than the others.
$ws.Columns("F").ColumnWidth = 22.42578125

# Page setup / view tidy-up to match the saved workbook state.
$ws.PageSetup.Orientation = 1
[void]$ws.Range("F18").Select()
